$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Cypher query for the "CasesTab" row (B2): the WITH clause now also
# carries demo.weight forward, the Age expression is no longer wrapped in a
# (redundant) coalesce, the Weight coalesce is re-indented, and the query now
# ends with an explicit ORDER BY / LIMIT clause.
$casesQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
WHERE diag.disease_term IN ['T Cell Lymphoma']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS ``Case ID``,
       coalesce(s.clinical_study_designation, '') AS ``Study Code``,
       coalesce(s.clinical_study_type, '') AS  ``Study Type``,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS ``Stage of Disease``,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS ``Weight (kg)``,
       coalesce(diag.best_response, '') AS ``Response to Treatment``,
       coalesce(co.cohort_description, '') AS ``Cohort``
Order by c.case_id LIMIT 100
"@

# Trim the trailing newline that the here-string literal adds after the
# closing line so the stored text ends exactly at "LIMIT 100".
$casesQuery = $casesQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $casesQuery

# The row grew by one wrapped line of text, so Excel would now auto-size it
# taller; pin the new height directly since this content is long & wrapped.
$ws.Rows.Item(2).RowHeight = 360

# Cursor/selection moved on to B5 in the saved file.
$ws.Range("B5").Select()
